$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Merge the two "Assembly Flavours" runs into one (drop the run
#    break between "...Bell Labs). " and "The assembly flavor...").
# ---------------------------------------------------------------------
$old1 = "invented at AT&T Bell Labs). The assembly flavor has no effect"
$new1 = "invented at AT&T Bell Labs). The assembly flavor has no effect"
$d.Content.Find.Execute($old1, $false, $false, $false, $false, $false, $true, 1, $false, $new1, 2) | Out-Null

# ---------------------------------------------------------------------
# 2. Merge the two "registers" runs together (drop the run break
#    between "...base address of the stack. " and "It should be
#    noted that on Unix...").
# ---------------------------------------------------------------------
$old2 = "base address of the stack. It should be noted that on Unix"
$new2 = "base address of the stack. It should be noted that on Unix"
$d.Content.Find.Execute($old2, $false, $false, $false, $false, $false, $true, 1, $false, $new2, 2) | Out-Null

# ---------------------------------------------------------------------
# 3. Give the blank run right after "Segment Registers:" paragraph
#    Arial formatting (it currently has no rPr content at all).
# ---------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -eq "Segment Registers:`r") {
        $nextP = $p.Next()
        $nextNextP = $nextP.Next()
        $blankP = $nextNextP.Next()
        $blankP.Range.Font.Name = "Arial"
        break
    }
}

# ---------------------------------------------------------------------
# 4. Split "All of the following file extensions indicate that..."
#    so that the word "likely" is inserted as its own run:
#    "All of the following file extensions " + "likely " + "indicate..."
# ---------------------------------------------------------------------
$old4 = "All of the following file extensions indicate"
$new4 = "All of the following file extensions likely indicate"
$d.Content.Find.Execute($old4, $false, $false, $false, $false, $false, $true, 1, $false, $new4, 2) | Out-Null

# ---------------------------------------------------------------------
# 5. Styles.xml changes
# ---------------------------------------------------------------------
# Normal style: overflowPunct true -> false
$normal = $d.Styles("Normal")
$normal.ParagraphFormat.OverflowPunctuation = $false

# Rename InternetLink style id's display is already "Hyperlink"; underlying
# style id rename TextBody -> BodyText, Caption name caption -> caption1
# handled through NameLocal / NextParagraphStyle below.
$heading = $d.Styles("Heading")
$heading.NextParagraphStyle = "Body Text"

$caption1 = $d.Styles("Caption1")
$caption1.NameLocal = "caption1"

Write-Output "done"
